$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "21년 2월 2주"
$ws.Range("E18").Value = "http://freesearch.pe.kr/archives/5034?utm_source=rss&utm_medium=rss&utm_campaign=21%eb%85%84-2%ec%9b%94-2%ec%a3%bc"

$ws.Range("D51").Value = "현재(2023-09-20) bskyvision이 티스토리 IT 분야 크리에이터 중에서 9번째로 소개되고 있습니다!"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%ED%98%84%EC%9E%AC2023-09-20-bskyvision%EC%9D%B4-%ED%8B%B0%EC%8A%A4%ED%86%A0%EB%A6%AC-IT-%EB%B6%84%EC%95%BC-%ED%81%AC%EB%A6%AC%EC%97%90%EC%9D%B4%ED%84%B0-%EC%A4%91%EC%97%90%EC%84%9C-9%EB%B2%88%EC%A7%B8%EB%A1%9C-%EC%86%8C%EA%B0%9C%EB%90%98%EA%B3%A0-%EC%9E%88%EC%8A%B5%EB%8B%88%EB%8B%A4"
